# Update Riolunato report: append rows for 29 June 2021 through 25 July 2021
# (serial dates 44376 .. 44402), mirroring the existing pattern of rows,
# with column A formatted like the preceding date cells and columns B:D
# holding zero values (aggiornamento fino a 28 luglio).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 301
$firstNewRow = 302
$lastNewRow = 328
$startSerial = 44376

# Copy the formatting (number format, borders, alignment, font) of the
# template row (the last existing data row) down across the new rows in
# one shot, then fill in the actual values.
$srcRow = $ws.Range("A" + $lastRow + ":D" + $lastRow)
$destRows = $ws.Range("A" + $firstNewRow + ":D" + $lastNewRow)
$srcRow.Copy()
$destRows.PasteSpecial(-4122)

$serial = $startSerial
for ($r = $firstNewRow; $r -le $lastNewRow; $r++) {
    $ws.Range("A" + $r).Value = $serial
    $ws.Range("B" + $r).Value = 0
    $ws.Range("C" + $r).Value = 0
    $ws.Range("D" + $r).Value = 0
    $serial = $serial + 1
}

$excel.CutCopyMode = 0

$addr = $ws.UsedRange.Address()
Write-Host ("UsedRange now: " + $addr)
